$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("sheet 1")

# Update the last row of sheet 1 with unicode test data.
$ws1.Range("A4").Value = "α"
$ws1.Range("A4").Font.Name = "Calibri"
$ws1.Range("B4").Value = "unicode"

# Make "sheet 1" the active sheet/tab and select A4 on it.
$ws1.Activate()
$ws1.Range("A4").Select()

$wb.Save()
